$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$r1 = $tr.Characters(1, 13)
$r1.Text = "Final CPU "
